$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mejoras")

# New improvement record in row 5: frmCombo file, with a new detailed note.
$newText = "' - Limitar la funcionalidad del botón ""Ingresar"", Ya que si ingresas el nombre de un combo que ya existe te permitirá agregar productos que dicho comobo no tiene.`n- Aún que en la aplicaciòn muestre un valor en la base de datos se recalcula el valor correcto del combo."

$ws.Range("C5").Value = "frmCombo"
$ws.Range("D5").Value = $newText

# Style D5: left/top aligned, wrap text (leading apostrophe above set the
# quote-prefix so Excel does not try to parse the leading " - " as a formula).
$ws.Range("D5").HorizontalAlignment = -4131  # xlLeft
$ws.Range("D5").VerticalAlignment = -4160    # xlTop
$ws.Range("D5").WrapText = $true

# Row height for row 5
$ws.Rows.Item(5).RowHeight = 14.25

# Column D width
$ws.Columns.Item(4).ColumnWidth = 39.14

# Freeze panes above row 2 (so header row stays visible) and set the
# final active selection to D5, matching the saved view state.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D5").Select()
